# Update the RPAR holdings workbook:
#  - bump the "as of" date in the confidential disclaimer (A18)
#  - refresh the Weight / Percent Change figures in D2:E15
#
# The sheet ships protected, so it must be unlocked before writing and
# re-locked afterwards to leave it in the same (protected) state it
# started in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.05687629715038708
$ws.Range("E2").Value = 0.005097076131782563

$ws.Range("D3").Value = 0.02357700543567781
$ws.Range("E3").Value = 0.008215962441314506

$ws.Range("D4").Value = 0.0311375391204085
$ws.Range("E4").Value = 0.009902875642734843

$ws.Range("D5").Value = 0.03266084664799868
$ws.Range("E5").Value = 0.007282483710233878

$ws.Range("D6").Value = 0.0380023060451326
$ws.Range("E6").Value = 0.02090918547799858

$ws.Range("D7").Value = 0.01926008894745511
$ws.Range("E7").Value = 0.01293582375478919

$ws.Range("D8").Value = 0.004336682589359249
$ws.Range("E8").Value = -0.01321786690975413

$ws.Range("D9").Value = 0.006926041838247406
$ws.Range("E9").Value = 0.006849315068493178

$ws.Range("D10").Value = 0.07016965903475539
$ws.Range("E10").Value = 0.01633802816901397

$ws.Range("D11").Value = 0.07020919123702851
$ws.Range("E11").Value = 0.01689189189189166

$ws.Range("D12").Value = 0.1472600889474551
$ws.Range("E12").Value = 0.001646502970863883

$ws.Range("D13").Value = 0.3844111349036403
$ws.Range("E13").Value = 0.0004370629370629153

$ws.Range("D14").Value = 0.1151731181024543
$ws.Range("E14").Value = -0.002522825564632347

$ws.Range("E15").Value = 0.004515987481469086

# Writing the multi-line disclaimer makes the host auto-mark row 18 with an
# explicit (custom) height; AutoFit puts it back to the sheet's implicit
# default so the row stays unchanged like every other data row.
$ws.Rows(18).AutoFit()

$ws.Protect()
